$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Update "Recommandations" sheet: rows 2-44 (BRVM data refresh) ---
$ws1.Cells.Item(2,1).Value = 'CFAO MOTORS CI'
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 3
$ws1.Cells.Item(2,4).Value = 2845
$ws1.Cells.Item(2,5).Value = 960
$ws1.Cells.Item(2,6).Value = '🟡 Observer'
$ws1.Cells.Item(2,7).Value = '➖ Neutre'

$ws1.Cells.Item(3,1).Value = 'UNIWAX CI'
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 3
$ws1.Cells.Item(3,4).Value = 2650
$ws1.Cells.Item(3,5).Value = 805
$ws1.Cells.Item(3,6).Value = '🟡 Observer'
$ws1.Cells.Item(3,7).Value = '➖ Neutre'

$ws1.Cells.Item(4,1).Value = 'BRVM - SERVICES PUBLICS'
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 6
$ws1.Cells.Item(4,4).Value = 2484.38
$ws1.Cells.Item(4,5).Value = 104
$ws1.Cells.Item(4,6).Value = '🟡 Observer'
$ws1.Cells.Item(4,7).Value = '➖ Neutre'

$ws1.Cells.Item(5,1).Value = 'SETAO CI'
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 3
$ws1.Cells.Item(5,4).Value = 2135
$ws1.Cells.Item(5,5).Value = 700
$ws1.Cells.Item(5,6).Value = '🟡 Observer'
$ws1.Cells.Item(5,7).Value = '➖ Neutre'

$ws1.Cells.Item(6,1).Value = 'NEI-CEDA CI'
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 1980
$ws1.Cells.Item(6,5).Value = 680
$ws1.Cells.Item(6,6).Value = '🟡 Observer'
$ws1.Cells.Item(6,7).Value = '➖ Neutre'

$ws1.Cells.Item(7,1).Value = 'BRVM - AUTRES SECTEURS'
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 1865.43
$ws1.Cells.Item(7,5).Value = 615.64
$ws1.Cells.Item(7,6).Value = '🟡 Observer'
$ws1.Cells.Item(7,7).Value = '➖ Neutre'

$ws1.Cells.Item(8,1).Value = 'AIR LIQUIDE CI'
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 3
$ws1.Cells.Item(8,4).Value = 1785
$ws1.Cells.Item(8,5).Value = 590
$ws1.Cells.Item(8,6).Value = '🟡 Observer'
$ws1.Cells.Item(8,7).Value = '➖ Neutre'

$ws1.Cells.Item(9,1).Value = 'BRVM - DISTRIBUTION'
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 3
$ws1.Cells.Item(9,4).Value = 1199.51
$ws1.Cells.Item(9,5).Value = 395.92
$ws1.Cells.Item(9,6).Value = '🟡 Observer'
$ws1.Cells.Item(9,7).Value = '➖ Neutre'

$ws1.Cells.Item(10,1).Value = 'BRVM - TRANSPORT'
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 3
$ws1.Cells.Item(10,4).Value = 1040.28
$ws1.Cells.Item(10,5).Value = 348.8
$ws1.Cells.Item(10,6).Value = '🟡 Observer'
$ws1.Cells.Item(10,7).Value = '➖ Neutre'

$ws1.Cells.Item(11,1).Value = 'BRVM - AGRICULTURE'
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 3
$ws1.Cells.Item(11,4).Value = 988.76
$ws1.Cells.Item(11,5).Value = 329.93
$ws1.Cells.Item(11,6).Value = '🟡 Observer'
$ws1.Cells.Item(11,7).Value = '➖ Neutre'

$ws1.Cells.Item(12,1).Value = 'Marché'
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 1
$ws1.Cells.Item(12,4).Value = 805
$ws1.Cells.Item(12,5).Value = 805
$ws1.Cells.Item(12,6).Value = '🟡 Observer'
$ws1.Cells.Item(12,7).Value = '➖ Neutre'

$ws1.Cells.Item(13,1).Value = 'BRVM - CONSOMMATION DE BASE          (**)'
$ws1.Cells.Item(13,2).Value = 0
$ws1.Cells.Item(13,3).Value = 2
$ws1.Cells.Item(13,4).Value = 419.67
$ws1.Cells.Item(13,5).Value = 210.16
$ws1.Cells.Item(13,6).Value = '🟡 Observer'
$ws1.Cells.Item(13,7).Value = '➖ Neutre'

$ws1.Cells.Item(14,1).Value = 'BRVM - INDUSTRIELS'
$ws1.Cells.Item(14,2).Value = 0
$ws1.Cells.Item(14,3).Value = 3
$ws1.Cells.Item(14,4).Value = 412.46
$ws1.Cells.Item(14,5).Value = 138.17
$ws1.Cells.Item(14,6).Value = '🟡 Observer'
$ws1.Cells.Item(14,7).Value = '➖ Neutre'

$ws1.Cells.Item(15,1).Value = 'BRVM-PRESTIGE'
$ws1.Cells.Item(15,2).Value = 0
$ws1.Cells.Item(15,3).Value = 3
$ws1.Cells.Item(15,4).Value = 400.35
$ws1.Cells.Item(15,5).Value = 134.13
$ws1.Cells.Item(15,6).Value = '🟡 Observer'
$ws1.Cells.Item(15,7).Value = '➖ Neutre'

$ws1.Cells.Item(16,1).Value = 'BRVM - FINANCES'
$ws1.Cells.Item(16,2).Value = 0
$ws1.Cells.Item(16,3).Value = 3
$ws1.Cells.Item(16,4).Value = 380.32
$ws1.Cells.Item(16,5).Value = 128.07
$ws1.Cells.Item(16,6).Value = '🟡 Observer'
$ws1.Cells.Item(16,7).Value = '➖ Neutre'

$ws1.Cells.Item(17,1).Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Cells.Item(17,2).Value = 0
$ws1.Cells.Item(17,3).Value = 3
$ws1.Cells.Item(17,4).Value = 374.54
$ws1.Cells.Item(17,5).Value = 122.66
$ws1.Cells.Item(17,6).Value = '🟡 Observer'
$ws1.Cells.Item(17,7).Value = '➖ Neutre'

$ws1.Cells.Item(18,1).Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 3
$ws1.Cells.Item(18,4).Value = 373.76
$ws1.Cells.Item(18,5).Value = 125.86
$ws1.Cells.Item(18,6).Value = '🟡 Observer'
$ws1.Cells.Item(18,7).Value = '➖ Neutre'

$ws1.Cells.Item(19,1).Value = 'BRVM - ENERGIE'
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 3
$ws1.Cells.Item(19,4).Value = 323.31
$ws1.Cells.Item(19,5).Value = 107.77
$ws1.Cells.Item(19,6).Value = '🟡 Observer'
$ws1.Cells.Item(19,7).Value = '➖ Neutre'

$ws1.Cells.Item(20,1).Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 3
$ws1.Cells.Item(20,4).Value = 283.47
$ws1.Cells.Item(20,5).Value = 94.42
$ws1.Cells.Item(20,6).Value = '🟡 Observer'
$ws1.Cells.Item(20,7).Value = '➖ Neutre'

$ws1.Cells.Item(21,1).Value = 'BRVM - INDUSTRIE                  (**)'
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = 251.95
$ws1.Cells.Item(21,5).Value = 251.95
$ws1.Cells.Item(21,6).Value = '🟡 Observer'
$ws1.Cells.Item(21,7).Value = '➖ Neutre'

$ws1.Cells.Item(22,1).Value = 'BRVM - INDUSTRIE             (**)'
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = 251.86
$ws1.Cells.Item(22,5).Value = 251.86
$ws1.Cells.Item(22,6).Value = '🟡 Observer'
$ws1.Cells.Item(22,7).Value = '➖ Neutre'

$ws1.Cells.Item(23,1).Value = 'BRVM-PRINCIPAL                (**)'
$ws1.Cells.Item(23,2).Value = 0
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = 191.84
$ws1.Cells.Item(23,5).Value = 191.84
$ws1.Cells.Item(23,6).Value = '🟡 Observer'
$ws1.Cells.Item(23,7).Value = '➖ Neutre'

$ws1.Cells.Item(24,1).Value = 'BRVM-PRINCIPAL                    (**)'
$ws1.Cells.Item(24,2).Value = 0
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = 191.6
$ws1.Cells.Item(24,5).Value = 191.6
$ws1.Cells.Item(24,6).Value = '🟡 Observer'
$ws1.Cells.Item(24,7).Value = '➖ Neutre'

$ws1.Cells.Item(25,1).Value = 'NEI-CEDA CI (NEIC)'
$ws1.Cells.Item(25,2).Value = 1
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 6.98
$ws1.Cells.Item(25,5).Value = 6.98
$ws1.Cells.Item(25,6).Value = '🟡 Observer'
$ws1.Cells.Item(25,7).Value = '➖ Neutre'

$ws1.Cells.Item(26,1).Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Cells.Item(26,2).Value = 1
$ws1.Cells.Item(26,3).Value = 0
$ws1.Cells.Item(26,4).Value = 6.93
$ws1.Cells.Item(26,5).Value = 6.93
$ws1.Cells.Item(26,6).Value = '🟡 Observer'
$ws1.Cells.Item(26,7).Value = '➖ Neutre'

$ws1.Cells.Item(27,1).Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Cells.Item(27,2).Value = 2
$ws1.Cells.Item(27,3).Value = 1
$ws1.Cells.Item(27,4).Value = 6.88
$ws1.Cells.Item(27,5).Value = -6.22
$ws1.Cells.Item(27,6).Value = '🟡 Observer'
$ws1.Cells.Item(27,7).Value = '👀 À surveiller'

$ws1.Cells.Item(28,1).Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 0
$ws1.Cells.Item(28,4).Value = 6.25
$ws1.Cells.Item(28,5).Value = 6.25
$ws1.Cells.Item(28,6).Value = '🟡 Observer'
$ws1.Cells.Item(28,7).Value = '➖ Neutre'

$ws1.Cells.Item(29,1).Value = 'BANK OF AFRICA BN (BOAB)'
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 4.01
$ws1.Cells.Item(29,5).Value = 4.01
$ws1.Cells.Item(29,6).Value = '🟡 Observer'
$ws1.Cells.Item(29,7).Value = '➖ Neutre'

$ws1.Cells.Item(30,1).Value = 'BICI CI (BICC)'
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 0
$ws1.Cells.Item(30,4).Value = 3.44
$ws1.Cells.Item(30,5).Value = 3.44
$ws1.Cells.Item(30,6).Value = '🟡 Observer'
$ws1.Cells.Item(30,7).Value = '➖ Neutre'

$ws1.Cells.Item(31,1).Value = 'NESTLE CI (NTLC)'
$ws1.Cells.Item(31,2).Value = 1
$ws1.Cells.Item(31,3).Value = 0
$ws1.Cells.Item(31,4).Value = 3.34
$ws1.Cells.Item(31,5).Value = 3.34
$ws1.Cells.Item(31,6).Value = '🟡 Observer'
$ws1.Cells.Item(31,7).Value = '➖ Neutre'

$ws1.Cells.Item(32,1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Cells.Item(32,2).Value = 1
$ws1.Cells.Item(32,3).Value = 0
$ws1.Cells.Item(32,4).Value = 3.23
$ws1.Cells.Item(32,5).Value = 3.23
$ws1.Cells.Item(32,6).Value = '🟡 Observer'
$ws1.Cells.Item(32,7).Value = '➖ Neutre'

$ws1.Cells.Item(33,1).Value = 'BANK OF AFRICA ML (BOAM)'
$ws1.Cells.Item(33,2).Value = 1
$ws1.Cells.Item(33,3).Value = 0
$ws1.Cells.Item(33,4).Value = 2.23
$ws1.Cells.Item(33,5).Value = 2.23
$ws1.Cells.Item(33,6).Value = '🟡 Observer'
$ws1.Cells.Item(33,7).Value = '➖ Neutre'

$ws1.Cells.Item(34,1).Value = 'SMB CI (SMBC)'
$ws1.Cells.Item(34,2).Value = 1
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = 0.19
$ws1.Cells.Item(34,5).Value = 7.07
$ws1.Cells.Item(34,6).Value = '🟡 Observer'
$ws1.Cells.Item(34,7).Value = '👀 À surveiller'

$ws1.Cells.Item(35,1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Cells.Item(35,2).Value = 1
$ws1.Cells.Item(35,3).Value = 1
$ws1.Cells.Item(35,4).Value = 0.09
$ws1.Cells.Item(35,5).Value = 3.14
$ws1.Cells.Item(35,6).Value = '🟡 Observer'
$ws1.Cells.Item(35,7).Value = '👀 À surveiller'

$ws1.Cells.Item(36,1).Value = 'TOTAL'
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 3
$ws1.Cells.Item(36,4).Value = 0
$ws1.Cells.Item(36,5).Value = 0
$ws1.Cells.Item(36,6).Value = '🟡 Observer'
$ws1.Cells.Item(36,7).Value = '➖ Neutre'

$ws1.Cells.Item(37,1).Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Cells.Item(37,2).Value = 1
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -0.9399999999999999
$ws1.Cells.Item(37,5).Value = 4
$ws1.Cells.Item(37,6).Value = '🟡 Observer'
$ws1.Cells.Item(37,7).Value = '👀 À surveiller'

$ws1.Cells.Item(38,1).Value = 'BERNABE CI (BNBC)'
$ws1.Cells.Item(38,2).Value = 0
$ws1.Cells.Item(38,3).Value = 1
$ws1.Cells.Item(38,4).Value = -1.96
$ws1.Cells.Item(38,5).Value = -1.96
$ws1.Cells.Item(38,6).Value = '🟡 Observer'
$ws1.Cells.Item(38,7).Value = '➖ Neutre'

$ws1.Cells.Item(39,1).Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws1.Cells.Item(39,2).Value = 0
$ws1.Cells.Item(39,3).Value = 1
$ws1.Cells.Item(39,4).Value = -2.34
$ws1.Cells.Item(39,5).Value = -2.34
$ws1.Cells.Item(39,6).Value = '🟡 Observer'
$ws1.Cells.Item(39,7).Value = '➖ Neutre'

$ws1.Cells.Item(40,1).Value = 'SETAO CI (STAC)'
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 2
$ws1.Cells.Item(40,4).Value = -4.85
$ws1.Cells.Item(40,5).Value = -2.11
$ws1.Cells.Item(40,6).Value = '🟡 Observer'
$ws1.Cells.Item(40,7).Value = '👀 À surveiller'

$ws1.Cells.Item(41,1).Value = 'ONATEL BF (ONTBF)'
$ws1.Cells.Item(41,2).Value = 1
$ws1.Cells.Item(41,3).Value = 2
$ws1.Cells.Item(41,4).Value = -5.77
$ws1.Cells.Item(41,5).Value = -6.04
$ws1.Cells.Item(41,6).Value = '🟡 Observer'
$ws1.Cells.Item(41,7).Value = '👀 À surveiller'

$ws1.Cells.Item(42,1).Value = 'SICOR CI (SICC)'
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 1
$ws1.Cells.Item(42,4).Value = -7.38
$ws1.Cells.Item(42,5).Value = -7.38
$ws1.Cells.Item(42,6).Value = '🟡 Observer'
$ws1.Cells.Item(42,7).Value = '➖ Neutre'

$ws1.Cells.Item(43,1).Value = 'UNIWAX CI (UNXC)'
$ws1.Cells.Item(43,2).Value = 0
$ws1.Cells.Item(43,3).Value = 2
$ws1.Cells.Item(43,4).Value = -13.92
$ws1.Cells.Item(43,5).Value = -7.47
$ws1.Cells.Item(43,6).Value = '🟡 Observer'
$ws1.Cells.Item(43,7).Value = '👀 À surveiller'

$ws1.Cells.Item(44,1).Value = 'UNILEVER CI (UNLC)'
$ws1.Cells.Item(44,2).Value = 0
$ws1.Cells.Item(44,3).Value = 2
$ws1.Cells.Item(44,4).Value = -14.16
$ws1.Cells.Item(44,5).Value = -6.67
$ws1.Cells.Item(44,6).Value = '🟡 Observer'
$ws1.Cells.Item(44,7).Value = '➖ Neutre'

# Remove now-obsolete rows 45-48 (used range shrinks from G48 to G44)
$ws1.Rows("45:48").Delete()

# --- Update "Top_YTD" sheet: column B values for rows 2-11 ---
$ws2.Cells.Item(2,2).Value = 474585.55
$ws2.Cells.Item(3,2).Value = 114989.5
$ws2.Cells.Item(4,2).Value = 94513.23
$ws2.Cells.Item(5,2).Value = 53360
$ws2.Cells.Item(6,2).Value = 43773.05
$ws2.Cells.Item(7,2).Value = 37503.47
$ws2.Cells.Item(8,2).Value = 33468.5
$ws2.Cells.Item(9,2).Value = 12386.71
$ws2.Cells.Item(10,2).Value = 8816.52
$ws2.Cells.Item(11,2).Value = 7827.79

Write-Output "Update complete"